$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update res_bus vm_pu results for the 380 kV case: columns B:F and I:N
# for data rows 2-25 (bus indices 0-23). Column G (slack, =1) and column A
# (index) are unchanged; column H has no data in this sheet.

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.02303035848436
$bf[0,2] = 1.042105295979751
$bf[0,3] = 1.034387117077619
$bf[0,4] = 1.047097730243341
$ws.Range("B2:F2").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.034113813033221
$in[0,1] = 1.028212867046043
$in[0,2] = 1.044882923899889
$in[0,3] = 1.03718673765476
$in[0,4] = 1.049861319745721
$in[0,5] = 1.013431309832355
$ws.Range("I2:N2").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.023934054689497
$bf[0,2] = 1.042651281205054
$bf[0,3] = 1.035155690111072
$bf[0,4] = 1.047870743835646
$ws.Range("B3:F3").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.034179311119238
$in[0,1] = 1.028755104314368
$in[0,2] = 1.045240367376724
$in[0,3] = 1.037764596103085
$in[0,4] = 1.050446209527859
$in[0,5] = 1.013614153680398
$ws.Range("I3:N3").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.024519712828117
$bf[0,2] = 1.043004512840821
$bf[0,3] = 1.035654042481793
$bf[0,4] = 1.048371639919545
$ws.Range("B4:F4").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.034220291378506
$in[0,1] = 1.029106261735243
$in[0,2] = 1.045470891952494
$in[0,3] = 1.038138882385756
$in[0,4] = 1.050824715417547
$in[0,5] = 1.013732476509294
$ws.Range("I4:N4").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.024766139084791
$bf[0,2] = 1.04315299556619
$bf[0,3] = 1.035863795599169
$bf[0,4] = 1.048582383190395
$ws.Range("B5:F5").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.034237183216239
$in[0,1] = 1.029253957417717
$in[0,2] = 1.045567619678836
$in[0,3] = 1.03829632025776
$in[0,4] = 1.050983847838367
$in[0,5] = 1.013782221497607
$ws.Range("I5:N5").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.024807527738854
$bf[0,2] = 1.043177925477267
$bf[0,3] = 1.03589902843045
$bf[0,4] = 1.048617777618027
$ws.Range("B6:F6").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.034239999689084
$in[0,1] = 1.029278760190116
$in[0,2] = 1.045583849805984
$in[0,3] = 1.038322759883177
$in[0,4] = 1.05101056732157
$in[0,5] = 1.013790574007684
$ws.Range("I6:N6").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.024523004742483
$bf[0,2] = 1.043006496939262
$bf[0,3] = 1.035656844249747
$bf[0,4] = 1.048374455229669
$ws.Range("B7:F7").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.034220518410321
$in[0,1] = 1.029108234982156
$in[0,2] = 1.045472185161291
$in[0,3] = 1.038140985734649
$in[0,4] = 1.050826841721402
$in[0,5] = 1.013733141196721
$ws.Range("I7:N7").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.023335578364427
$bf[0,2] = 1.042289824829538
$bf[0,3] = 1.034646644698654
$bf[0,4] = 1.047358826636917
$ws.Range("B8:F8").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.034136238121153
$in[0,1] = 1.028396057147096
$in[0,2] = 1.045003881085834
$in[0,3] = 1.037381949482795
$in[0,4] = 1.050058975999698
$in[0,5] = 1.013493100280031
$ws.Range("I8:N8").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.021250183753032
$bf[0,2] = 1.041026612168894
$bf[0,3] = 1.032874545193919
$bf[0,4] = 1.045574649608362
$ws.Range("B9:F9").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033977026547863
$in[0,1] = 1.027143410803385
$in[0,2] = 1.044172879212099
$in[0,3] = 1.036047352429641
$in[0,4] = 1.048706306764949
$in[0,5] = 1.013070219926541
$ws.Range("I9:N9").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.019864711170086
$bf[0,2] = 1.040184365703948
$bf[0,3] = 1.031698633340708
$bf[0,4] = 1.044389016434975
$ws.Range("B10:F10").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033863737984512
$in[0,1] = 1.026309932324711
$in[0,2] = 1.043615080349202
$in[0,3] = 1.035159670916014
$in[0,4] = 1.04780490739023
$in[0,5] = 1.012788396813636
$ws.Range("I10:N10").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.019265939750878
$bf[0,2] = 1.039819665770582
$bf[0,3] = 1.031190775039438
$bf[0,4] = 1.043876556612105
$ws.Range("B11:F11").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033812997156808
$in[0,1] = 1.025949425986423
$in[0,2] = 1.043372666708751
$in[0,3] = 1.034775799236199
$in[0,4] = 1.047414703980174
$in[0,5] = 1.012666393287281
$ws.Range("I11:N11").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.019043703137354
$bf[0,2] = 1.039684201911456
$bf[0,3] = 1.031002333918701
$bf[0,4] = 1.043686347686932
$ws.Range("B12:F12").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033793897195289
$in[0,1] = 1.025815578482735
$in[0,2] = 1.043282492444173
$in[0,3] = 1.034633288816565
$in[0,4] = 1.0472697831575
$in[0,5] = 1.012621080391725
$ws.Range("I12:N12").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.019091365760575
$bf[0,2] = 1.03971325924657
$bf[0,3] = 1.031042746095638
$bf[0,4] = 1.043727141710355
$ws.Range("B13:F13").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033798005618448
$in[0,1] = 1.025844286464895
$in[0,2] = 1.043301841034074
$in[0,3] = 1.034663854295053
$in[0,4] = 1.047300868318053
$in[0,5] = 1.012630799943701
$ws.Range("I13:N13").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.019247566050617
$bf[0,2] = 1.039808468233671
$bf[0,3] = 1.031175194345625
$bf[0,4] = 1.043860830984501
$ws.Range("B14:F14").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033811423495745
$in[0,1] = 1.025938360866776
$in[0,2] = 1.043365215539819
$in[0,3] = 1.034764017709147
$in[0,4] = 1.047402724400169
$in[0,5] = 1.012662647611648
$ws.Range("I14:N14").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.01934382926962
$bf[0,2] = 1.039867129954412
$bf[0,3] = 1.031256826649802
$bf[0,4] = 1.043943220169773
$ws.Range("B15:F15").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033819657249106
$in[0,1] = 1.025996331277054
$in[0,2] = 1.043404245348172
$in[0,3] = 1.034825741889969
$in[0,4] = 1.04746548374599
$in[0,5] = 1.012682270640972
$ws.Range("I15:N15").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.019904473947808
$bf[0,2] = 1.040208569789051
$bf[0,3] = 1.031732366177989
$bf[0,4] = 1.04442304642564
$ws.Range("B16:F16").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033867070026461
$in[0,1] = 1.026333866375434
$in[0,2] = 1.043631150077946
$in[0,3] = 1.035185157879436
$in[0,4] = 1.047830806344578
$in[0,5] = 1.012796494411198
$ws.Range("I16:N16").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.020256459690035
$bf[0,2] = 1.040422747076768
$bf[0,3] = 1.032031014082845
$bf[0,4] = 1.044724278523407
$ws.Range("B17:F17").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033896359922902
$in[0,1] = 1.026545699920749
$in[0,2] = 1.043773246149656
$in[0,3] = 1.035410744997753
$in[0,4] = 1.048059993871424
$in[0,5] = 1.012868151729479
$ws.Range("I17:N17").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.020461877628617
$bf[0,2] = 1.040547672759287
$bf[0,3] = 1.032205337521052
$bf[0,4] = 1.044900071421927
$ws.Range("B18:F18").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033913281485505
$in[0,1] = 1.026669296907722
$in[0,2] = 1.043856043073853
$in[0,3] = 1.035542374318886
$in[0,4] = 1.048193685423122
$in[0,5] = 1.012909950875432
$ws.Range("I18:N18").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.020531938556871
$bf[0,2] = 1.040590269070262
$bf[0,3] = 1.032264798821977
$bf[0,4] = 1.044960027340161
$ws.Range("B19:F19").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033919023674555
$in[0,1] = 1.026711446705876
$in[0,2] = 1.04388426016577
$in[0,3] = 1.035587264642304
$in[0,4] = 1.048239272487175
$in[0,5] = 1.012924203738017
$ws.Range("I19:N19").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.020218683500183
$bf[0,2] = 1.040399767896619
$bf[0,3] = 1.031998958822047
$bf[0,4] = 1.044691949916461
$ws.Range("B20:F20").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033893234221501
$in[0,1] = 1.026522968240271
$in[0,2] = 1.043758009398272
$in[0,3] = 1.035386536641779
$in[0,4] = 1.048035403135781
$in[0,5] = 1.012860463300394
$ws.Range("I20:N20").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.019201564150913
$bf[0,2] = 1.039780431489096
$bf[0,3] = 1.031136186099522
$bf[0,4] = 1.043821458902664
$ws.Range("B21:F21").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033807479231048
$in[0,1] = 1.025910656617564
$in[0,2] = 1.043346556931914
$in[0,3] = 1.034734519953605
$in[0,4] = 1.047372729808642
$in[0,5] = 1.012653269135116
$ws.Range("I21:N21").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.018563066806547
$bf[0,2] = 1.039391042375298
$bf[0,3] = 1.030594885076326
$bf[0,4] = 1.043274966230158
$ws.Range("B22:F22").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033752100814226
$in[0,1] = 1.025526023219452
$in[0,2] = 1.043087103049904
$in[0,3] = 1.034325015168299
$in[0,4] = 1.046956186383701
$in[0,5] = 1.012523025058976
$ws.Range("I22:N22").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.018901450515405
$bf[0,2] = 1.039597463022225
$bf[0,3] = 1.030881728565155
$bf[0,4] = 1.043564593924192
$ws.Range("B23:F23").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033781596177944
$in[0,1] = 1.025729890966548
$in[0,2] = 1.043224715712319
$in[0,3] = 1.03454205883314
$in[0,4] = 1.047176993324215
$in[0,5] = 1.012592067177703
$ws.Range("I23:N23").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.020235752592701
$bf[0,2] = 1.040410151200222
$bf[0,3] = 1.032013442821903
$bf[0,4] = 1.044706557546681
$ws.Range("B24:F24").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.033894647094345
$in[0,1] = 1.026533239590736
$in[0,2] = 1.043764894494156
$in[0,3] = 1.03539747520714
$in[0,4] = 1.048046514598977
$in[0,5] = 1.01286393736244
$ws.Range("I24:N24").Value = $in

$bf = New-Object 'object[,]' 1,5
$bf[0,0] = 1.02
$bf[0,1] = 1.02178846995074
$bf[0,2] = 1.04135321065727
$bf[0,3] = 1.03333171605175
$bf[0,4] = 1.046035238430111
$ws.Range("B25:F25").Value = $bf

$in = New-Object 'object[,]' 1,6
$in[0,0] = 1.034019449706743
$in[0,1] = 1.027466969762562
$in[0,2] = 1.044388389442056
$in[0,3] = 1.036392022708561
$in[0,4] = 1.049055944911968
$in[0,5] = 1.013179529585157
$ws.Range("I25:N25").Value = $in
